$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force a literal text value, bypassing Excel's automatic
    # number/date detection, and leave the cell's style index
    # untouched (reset back to the default "Normal" style after).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 52
Set-TextCell $ws.Range("A52") "2024-04-28"
$ws.Range("B52").Value = "Cteep"
$ws.Range("C52").Value = "a"
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 1
$ws.Range("F52").Value = 1
Set-TextCell $ws.Range("G52") "1"
$ws.Range("H52").Value = "Acessos"
$ws.Range("I52").Value = 1

# Row 53
Set-TextCell $ws.Range("A53") "2024-04-28"
$ws.Range("B53").Value = "Cteep"
$ws.Range("C53").Value = "asd"
$ws.Range("D53").Value = "asd"
$ws.Range("E53").Value = 123
$ws.Range("F53").Value = 123
Set-TextCell $ws.Range("G53") "123"
$ws.Range("H53").Value = "Acessos"
$ws.Range("I53").Value = 123
$ws.Range("J53").Value = ""

# Row 54
Set-TextCell $ws.Range("A54") "2024-04-22"
$ws.Range("B54").Value = "Henkel"
Set-TextCell $ws.Range("C54") "333"
Set-TextCell $ws.Range("D54") "333"
$ws.Range("E54").Value = 333
$ws.Range("F54").Value = 333
Set-TextCell $ws.Range("G54") "333"
$ws.Range("H54").Value = "Acessos"
$ws.Range("I54").Value = 333
Set-TextCell $ws.Range("J54") "333"

# Row 55
Set-TextCell $ws.Range("A55") "2024-04-29"
$ws.Range("B55").Value = "Henkel"
Set-TextCell $ws.Range("C55") "123"
Set-TextCell $ws.Range("D55") "123"
$ws.Range("E55").Value = 123
$ws.Range("F55").Value = 123
Set-TextCell $ws.Range("G55") "123"
$ws.Range("H55").Value = "Acessos"
$ws.Range("I55").Value = 123
$ws.Range("J55").Value = ""

# Row 56
Set-TextCell $ws.Range("A56") "2024-04-29"
$ws.Range("B56").Value = "Cteep"
Set-TextCell $ws.Range("C56") "123"
Set-TextCell $ws.Range("D56") "123"
$ws.Range("E56").Value = 123
$ws.Range("F56").Value = 123
Set-TextCell $ws.Range("G56") "123"
$ws.Range("H56").Value = "Acessos"
$ws.Range("I56").Value = 123
$ws.Range("J56").Value = ""

# Row 57
Set-TextCell $ws.Range("A57") "2024-04-29"
$ws.Range("B57").Value = "Cteep"
Set-TextCell $ws.Range("C57") "123"
Set-TextCell $ws.Range("D57") "123"
$ws.Range("E57").Value = 123
$ws.Range("F57").Value = 123
Set-TextCell $ws.Range("G57") "13"
$ws.Range("H57").Value = "Acessos"
$ws.Range("I57").Value = 123
$ws.Range("J57").Value = ""

# Row 58
Set-TextCell $ws.Range("A58") "2024-04-29"
$ws.Range("B58").Value = "Cteep"
Set-TextCell $ws.Range("C58") "123"
Set-TextCell $ws.Range("D58") "123"
$ws.Range("E58").Value = 123
$ws.Range("F58").Value = 123
Set-TextCell $ws.Range("G58") "123"
$ws.Range("H58").Value = "Acessos"
$ws.Range("I58").Value = 123
Set-TextCell $ws.Range("J58") "123"

# Row 59
Set-TextCell $ws.Range("A59") "2024-04-29"
$ws.Range("B59").Value = "Cteep"
Set-TextCell $ws.Range("C59") "123"
Set-TextCell $ws.Range("D59") "123"
$ws.Range("E59").Value = 123
$ws.Range("F59").Value = 123
Set-TextCell $ws.Range("G59") "123"
$ws.Range("H59").Value = "Acessos"
$ws.Range("I59").Value = 123
Set-TextCell $ws.Range("J59") "13"

# Row 60
Set-TextCell $ws.Range("A60") "2024-04-29"
$ws.Range("B60").Value = "Cteep"
Set-TextCell $ws.Range("C60") "123"
Set-TextCell $ws.Range("D60") "123"
$ws.Range("E60").Value = 123
$ws.Range("F60").Value = 123
Set-TextCell $ws.Range("G60") "123"
$ws.Range("H60").Value = "Acessos"
$ws.Range("I60").Value = 12
$ws.Range("J60").Value = ""

# Row 61
Set-TextCell $ws.Range("A61") "2024-04-29"
$ws.Range("B61").Value = "Flowserve"
Set-TextCell $ws.Range("C61") "234"
Set-TextCell $ws.Range("D61") "234"
$ws.Range("E61").Value = 234
$ws.Range("F61").Value = 234
Set-TextCell $ws.Range("G61") "234"
$ws.Range("H61").Value = "Acessos"
$ws.Range("I61").Value = 23
$ws.Range("J61").Value = ""

# Row 62
Set-TextCell $ws.Range("A62") "2024-04-29"
$ws.Range("B62").Value = "Flowserve"
Set-TextCell $ws.Range("C62") "123"
Set-TextCell $ws.Range("D62") "123"
$ws.Range("E62").Value = 123
$ws.Range("F62").Value = 123
Set-TextCell $ws.Range("G62") "13"
$ws.Range("H62").Value = "Acessos"
$ws.Range("I62").Value = 123
$ws.Range("J62").Value = ""

# Row 63
Set-TextCell $ws.Range("A63") "2024-04-29"
$ws.Range("B63").Value = "Flowserve"
Set-TextCell $ws.Range("C63") "123"
Set-TextCell $ws.Range("D63") "123"
$ws.Range("E63").Value = 123
$ws.Range("F63").Value = 132
Set-TextCell $ws.Range("G63") "123"
$ws.Range("H63").Value = "Acessos"
$ws.Range("I63").Value = 123
$ws.Range("J63").Value = ""

# Row 64
Set-TextCell $ws.Range("A64") "2024-04-29"
$ws.Range("B64").Value = "Flowserve"
Set-TextCell $ws.Range("C64") "123"
Set-TextCell $ws.Range("D64") "123"
$ws.Range("E64").Value = 123
$ws.Range("F64").Value = 123
Set-TextCell $ws.Range("G64") "123"
$ws.Range("H64").Value = "Acessos"
$ws.Range("I64").Value = 123
$ws.Range("J64").Value = ""

# Row 65
Set-TextCell $ws.Range("A65") "2024-04-29"
$ws.Range("B65").Value = "Henkel"
Set-TextCell $ws.Range("C65") "123"
Set-TextCell $ws.Range("D65") "123"
$ws.Range("E65").Value = 123
$ws.Range("F65").Value = 13
Set-TextCell $ws.Range("G65") "123"
$ws.Range("H65").Value = "Acessos"
$ws.Range("I65").Value = 13
$ws.Range("J65").Value = ""

# Row 66
Set-TextCell $ws.Range("A66") "2024-04-29"
$ws.Range("B66").Value = "Cteep"
Set-TextCell $ws.Range("C66") "123"
Set-TextCell $ws.Range("D66") "123"
$ws.Range("E66").Value = 123
$ws.Range("F66").Value = 123
Set-TextCell $ws.Range("G66") "123"
$ws.Range("H66").Value = "Acessos"
$ws.Range("I66").Value = 123
$ws.Range("J66").Value = ""

# Row 67
Set-TextCell $ws.Range("A67") "2024-04-29"
$ws.Range("B67").Value = "Cteep"
$ws.Range("C67").Value = "wer"
Set-TextCell $ws.Range("D67") "234"
$ws.Range("E67").Value = 234
$ws.Range("F67").Value = 234
Set-TextCell $ws.Range("G67") "234"
$ws.Range("H67").Value = "Acessos"
$ws.Range("I67").Value = 234
$ws.Range("J67").Value = "234`n234`n234`n"

# Row 68
Set-TextCell $ws.Range("A68") "2024-04-29"
$ws.Range("B68").Value = "Flowserve"
Set-TextCell $ws.Range("C68") "123"
Set-TextCell $ws.Range("D68") "123"
$ws.Range("E68").Value = 123
$ws.Range("F68").Value = 123
Set-TextCell $ws.Range("G68") "123"
$ws.Range("H68").Value = "Acessos"
$ws.Range("I68").Value = 123
$ws.Range("J68").Value = ""

# Row 69
Set-TextCell $ws.Range("A69") "2024-04-29"
$ws.Range("B69").Value = "Cteep"
Set-TextCell $ws.Range("C69") "123"
Set-TextCell $ws.Range("D69") "123"
Set-TextCell $ws.Range("E69") "123"
Set-TextCell $ws.Range("F69") "123"
Set-TextCell $ws.Range("G69") "123"
$ws.Range("H69").Value = "Acessos"
Set-TextCell $ws.Range("I69") "123"
Set-TextCell $ws.Range("J69") "123"

